# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G (header "K") with recalculated strikeout counts for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 2
    15 = 1
    16 = 3
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 2
    23 = 1
    24 = 4
    25 = 2
    26 = 1
    27 = 0
    28 = 0
    29 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
